$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 176 (01-07-2021) with revised figures ---
$ws.Cells.Item(176, 2).Value = 378
$ws.Cells.Item(176, 3).Value = 178
$ws.Cells.Item(176, 4).Value = 209
$ws.Cells.Item(176, 5).Value = 99

# --- Add new row 177 for 01-08-2021 ---
# Column A holds a date-formatted label ("01-08-2021") that must be stored
# as plain text (a shared string), not auto-converted to a date serial
# number. Typing the literal text directly into a cell causes Excel to
# recognize the pattern and convert it to a date. To avoid that (and avoid
# introducing any new/unused cell style), compute the text via a formula
# in a scratch cell (formula results are not subject to that date
# auto-conversion), then copy/paste-special as a value into the target
# cell, and finally remove the scratch cell.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.Formula = "=""01-08-2021"""
$scratch.Copy()
$ws.Cells.Item(177, 1).PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false

$ws.Cells.Item(177, 2).Value = 368
$ws.Cells.Item(177, 3).Value = 170
$ws.Cells.Item(177, 4).Value = 206
$ws.Cells.Item(177, 5).Value = 95
